$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp text (row 1)
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 22:52"

# Update Estados Unidos (row 4) statistics
$ws.Range("B4").Value = 609614
$ws.Range("C4").Value = 22673
$ws.Range("E4").Value = 545341
$ws.Range("F4").Value = 13438
$ws.Range("G4").Value = 2154
$ws.Range("H4").Value = 25794

# Update Brasil (row 17) statistics
$ws.Range("B17").Value = 24920
$ws.Range("C17").Value = 1490
$ws.Range("E17").Value = 20385
$ws.Range("G17").Value = 161
$ws.Range("H17").Value = 1489

# Reorder Letonia / Principado de Andorra in the country list and refresh
# Andorra's stats (row 88 keeps the Andorra data, which is now listed before
# Letonia; row 89 keeps Letonia's data, which did not change)
$ws.Range("A88").Value = "Principado de Andorra"
$ws.Range("B88").Value = 659
$ws.Range("C88").Value = 13
$ws.Range("D88").Value = 128
$ws.Range("E88").Value = 500
$ws.Range("F88").Value = 17
$ws.Range("G88").Value = 2
$ws.Range("H88").Value = 31

$ws.Range("A89").Value = "Letonia"
$ws.Range("B89").Value = 657
$ws.Range("C89").Value = 2
$ws.Range("D89").Value = 16
$ws.Range("E89").Value = 636
$ws.Range("F89").Value = 3
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 5
